$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 35818.75
$ws.Range("J3").Value = 35818.75
$ws.Range("L3").Value = 35818.75
$ws.Range("N3").Value = -36046.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 382.55554
$ws.Range("I6").Value = 76.15385000000001
$ws.Range("K6").Value = 228.46155
$ws.Range("M6").Value = -116.46155

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3134.6667
$ws.Range("J29").Value = 4552
$ws.Range("L29").Value = 13656
$ws.Range("N29").Value = -14218

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 68666.5
$ws.Range("J75").Value = 68666.5
$ws.Range("L75").Value = 68666.5
$ws.Range("N75").Value = -70538.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 68666.5
$ws.Range("J78").Value = 68666.5
$ws.Range("L78").Value = 205999.5
$ws.Range("N78").Value = -215359.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H102").Value = 35818.75
$ws.Range("J102").Value = 35818.75
$ws.Range("L102").Value = 35818.75
$ws.Range("N102").Value = -42308.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 7822.0713
$ws.Range("I131").Value = 3834.7778
$ws.Range("K131").Value = 11504.3334
$ws.Range("M131").Value = -6464.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3760.4482
$ws.Range("J138").Value = 3802.2307
$ws.Range("L138").Value = 11406.6921
$ws.Range("N138").Value = -21686.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 22729608
$ws.Range("I45").Value = 33334432
$ws.Range("J45").Value = 4989.4287
$ws.Range("K45").Value = 33334432
$ws.Range("L45").Value = 4989.4287
$ws.Range("M45").Value = -33334055
$ws.Range("N45").Value = -5743.4287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3461.3333
$ws.Range("I74").Value = 2198.077
$ws.Range("J74").Value = 5514.125
$ws.Range("K74").Value = 2198.077
$ws.Range("L74").Value = 5514.125
$ws.Range("M74").Value = -1324.077
$ws.Range("N74").Value = -7262.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3461.3333
$ws.Range("I77").Value = 2198.077
$ws.Range("J77").Value = 5514.125
$ws.Range("K77").Value = 10990.385
$ws.Range("L77").Value = 27570.625
$ws.Range("M77").Value = -6622.385000000002
$ws.Range("N77").Value = -36306.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 10568.125
$ws.Range("I110").Value = 10007.077
$ws.Range("K110").Value = 10007.077
$ws.Range("M110").Value = -7962.076999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 16667591
$ws.Range("I122").Value = 792.75
$ws.Range("K122").Value = 2378.25
$ws.Range("M122").Value = 71.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 50073.5
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 50073.5
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 50073.5
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -50693.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 69000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 69000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 142859490
$ws.Range("I86").Value = 500000740
$ws.Range("K86").Value = 500000740
$ws.Range("M86").Value = -499999617

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 142859490
$ws.Range("I89").Value = 500000740
$ws.Range("K89").Value = 2500003700
$ws.Range("M89").Value = -2499998084

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2001.4138
$ws.Range("I107").Value = 1838.7
$ws.Range("J107").Value = 2363
$ws.Range("K107").Value = 1838.7
$ws.Range("L107").Value = 2363
$ws.Range("M107").Value = 81.29999999999995
$ws.Range("N107").Value = -6203

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 5481.9165
$ws.Range("J116").Value = 4950
$ws.Range("L116").Value = 14850
$ws.Range("N116").Value = -21734

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 56789
$ws.Range("I123").Value = 56789
$ws.Range("K123").Value = 170367
$ws.Range("M123").Value = -167917

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 36797.668
$ws.Range("J136").Value = 36797.668
$ws.Range("L136").Value = 110393.004
$ws.Range("N136").Value = -115493.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 35716704
$ws.Range("I7").Value = 45456950
$ws.Range("J7").Value = 2469
$ws.Range("K7").Value = 45456950
$ws.Range("L7").Value = 2469
$ws.Range("M7").Value = -45456838
$ws.Range("N7").Value = -2693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1048.909
$ws.Range("I16").Value = 1220.8572
$ws.Range("J16").Value = 748
$ws.Range("K16").Value = 1220.8572
$ws.Range("L16").Value = 748
$ws.Range("M16").Value = -1050.8572
$ws.Range("N16").Value = -1088

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2800.4333
$ws.Range("I40").Value = 2458.72
$ws.Range("J40").Value = 4509
$ws.Range("K40").Value = 2458.72
$ws.Range("L40").Value = 4509
$ws.Range("M40").Value = -2322.72
$ws.Range("N40").Value = -4781

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3108.85
$ws.Range("J46").Value = 4949.727
$ws.Range("L46").Value = 4949.727
$ws.Range("N46").Value = -5325.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 55569250
$ws.Range("I61").Value = 62514812
$ws.Range("J61").Value = 4737.5
$ws.Range("K61").Value = 62514812
$ws.Range("L61").Value = 4737.5
$ws.Range("M61").Value = -62514610
$ws.Range("N61").Value = -5141.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5846.357
$ws.Range("I100").Value = 4234.8
$ws.Range("K100").Value = 4234.8
$ws.Range("M100").Value = -3693.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 55569250
$ws.Range("I113").Value = 62514812
$ws.Range("J113").Value = 4737.5
$ws.Range("K113").Value = 62514812
$ws.Range("L113").Value = 4737.5
$ws.Range("M113").Value = -62512642
$ws.Range("N113").Value = -9077.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3731
$ws.Range("I122").Value = 2846.5293
$ws.Range("K122").Value = 8539.5879
$ws.Range("M122").Value = -6089.5879

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 35716704
$ws.Range("I126").Value = 45456950
$ws.Range("J126").Value = 2469
$ws.Range("K126").Value = 136370850
$ws.Range("L126").Value = 7407
$ws.Range("M126").Value = -136368380
$ws.Range("N126").Value = -12347

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 259358.67
$ws.Range("J54").Value = 259358.67
$ws.Range("L54").Value = 259358.67
$ws.Range("N54").Value = -260398.67

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 214309280
$ws.Range("I107").Value = 500002500
$ws.Range("J107").Value = 100032000
$ws.Range("K107").Value = 1500007500
$ws.Range("L107").Value = 300096000
$ws.Range("M107").Value = -1500005580
$ws.Range("N107").Value = -300099840
